$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C: coin name & link swaps (rows 43/44 and 47/48) ---
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"

# --- Column D: price updates ---
# Rows where the new value would be misread as a pure number must be
# forced to Text (format "@") first so they round-trip as strings, just
# like the original inline-string cells.
$numericPriceRows = @(4,5,6,8,10,12,13,15,19,20,21,22,23,24,29,30,32,33,35,36,37,38,41,43,44,45,46,47,48,49,50,51)
foreach ($r in $numericPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "69.350.03"
$ws.Range("D3").Value = "3.683.83"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "685.71"
$ws.Range("D6").Value = "159.62"
$ws.Range("D7").Value = "3.681.74"
$ws.Range("D8").Value = "0.999"
$ws.Range("D10").Value = "0.146"
$ws.Range("D12").Value = "0.435"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D14").Value = "4.302.21"
$ws.Range("D15").Value = "32.41"
$ws.Range("D16").Value = "3.679.15"
$ws.Range("D17").Value = "69.329.73"
$ws.Range("D19").Value = "15.90"
$ws.Range("D20").Value = "6.45"
$ws.Range("D21").Value = "467.43"
$ws.Range("D22").Value = "10.13"
$ws.Range("D23").Value = "0.646"
$ws.Range("D24").Value = "79.33"
$ws.Range("D25").Value = "3.828.60"
$ws.Range("D29").Value = "9.18"
$ws.Range("D30").Value = "2.70"
$ws.Range("D32").Value = "6.62"
$ws.Range("D33").Value = "2.01"
$ws.Range("D35").Value = "26.71"
$ws.Range("D36").Value = "0.162"
$ws.Range("D37").Value = "8.17"
$ws.Range("D38").Value = "6.15"
$ws.Range("D41").Value = "0.0903"
$ws.Range("D43").Value = "0.942"
$ws.Range("D44").Value = "165.98"
$ws.Range("D45").Value = "47.70"
$ws.Range("D46").Value = "2.73"
$ws.Range("D47").Value = "1.12"
$ws.Range("D48").Value = "1.31"
$ws.Range("D49").Value = "0.000275"
$ws.Range("D50").Value = "28.10"
$ws.Range("D51").Value = "7.80"

# --- Column E: 1h volume/change percentages ---
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("E6").Value = "  -6.06%  "
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -5.61%  "
$ws.Range("E10").Value = "  -8.83%  "
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("E12").Value = "  -9.05%  "
$ws.Range("E13").Value = "  -6.68%  "
$ws.Range("E15").Value = "  -10.40%  "
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  -9.31%  "
$ws.Range("E20").Value = "  -10.10%  "
$ws.Range("E21").Value = "  -9.27%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -9.39%  "
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -11.43%  "
$ws.Range("E28").Value = "  -12.85%  "
$ws.Range("E29").Value = "  -10.11%  "
$ws.Range("E30").Value = "  -8.42%  "
$ws.Range("E31").Value = "  -12.81%  "
$ws.Range("E32").Value = "  -9.06%  "
$ws.Range("E33").Value = "  -10.79%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -8.38%  "
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("E37").Value = "  -12.04%  "
$ws.Range("E38").Value = "  -6.21%  "
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E41").Value = "  -9.84%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("E46").Value = "  -14.40%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("E49").Value = "  -8.70%  "
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("E51").Value = "  -9.51%  "
